$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Memory Usage (bytes)"

$ws.Range("C2").Value = 17.73285865783691
$ws.Range("C3").Value = 16.88194274902344
$ws.Range("C4").Value = 17.92383193969727
$ws.Range("C5").Value = 15.71798324584961
$ws.Range("C6").Value = 15.77877998352051
